# Refresh the "cryptos" price/volume snapshot (GitHub Actions data pull).
# Price cells (col D) are stored as plain text so values like "41.755.74"
# or "0.999" round-trip verbatim instead of being re-parsed as numbers;
# NumberFormat "@" forces text entry and Style "Normal" puts the cell back
# on the default (unstyled) format afterwards. Volume cells (col E) are
# already non-numeric ("  +0.54%  ") so a plain .Value assignment is safe.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.755.74"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.467.21"
$ws.Range("E3").Value = "  -0.91%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.78"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "92.86"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  +0.13%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.514"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.21%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "32.77"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.22%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0843"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.70%  "
$ws.Range("D13").Value = "2.849.38"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("E14").Value = "  +0.88%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.80"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").Value = "2.477.84"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").Value = "41.719.80"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("E20").Value = "  +3.08%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.68"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.00%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "71.20"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.75%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "239.54"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  -0.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.83"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("E29").Value = "  +1.35%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "35.91"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.36%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "155.98"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("E33").Value = "  +0.15%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0765"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.22%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.51"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.47%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.64"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("E40").Value = "  -2.34%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.02"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "1.976.61"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0284"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.94"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.37%  "
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "2.703.03"
$ws.Range("E48").Value = "  -0.85%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "96.97"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "67.22"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "73.17"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
